$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 (Model Evaluation): status flips from "On Going" to "Done", Kendala filled in ---
$ws.Range("G7").Value = "Done"
$ws.Range("H7").Value = "Tidak ada"

# --- Row 8 (classification report data): fill remaining columns ---
$ws.Range("D8").Value = "28 Desember 2021"
$ws.Range("F8").Value = "Muhammad Risky Pratama Hermawan"
$ws.Range("G8").Value = "Done"
$ws.Range("H8").Value = "Tidak ada"

# --- Row 9 (Deployment): fill remaining columns ---
# C9 and E9 are brand-new cells that pick up the new "text date, thin left/right
# border" style used elsewhere on this new row.
$ws.Range("C9").Value = "29 Desember 2021"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Borders(7).LineStyle = 1
$ws.Range("C9").Borders(10).LineStyle = 1
$ws.Range("C9").HorizontalAlignment = -4108
$ws.Range("C9").VerticalAlignment = -4108

$ws.Range("E9").Value = "2 Januari 2022"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Borders(7).LineStyle = 1
$ws.Range("E9").Borders(10).LineStyle = 1
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").VerticalAlignment = -4108

$ws.Range("D9").Value = "1 Januari 2022"

$ws.Range("F9").Value = "Muhammad Risky Pratama Hermawan"
$ws.Range("G9").Value = "Done"
$ws.Range("H9").Value = "Tidak ada"

# --- Row 10: stray Penanggung Jawab value from row 9 is cleared ---
$ws.Range("F10").ClearContents()

# --- Sheet view: scrolled one column right, selection moved to H10 ---
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("H10").Select()
